$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 771, shifting existing rows (771-812) down to (772-813)
$ws.Rows.Item(771).Insert()

# Ensure the date column stays plain text (matches the original inline-string
# storage) instead of being auto-converted into a date serial number.
$ws.Range("A771").NumberFormat = "@"

# Populate the newly inserted row with the new record
$ws.Range("A771").Value = "2026/02/08"
$ws.Range("B771").Value = "日"
$ws.Range("C771").Value = 0
$ws.Range("D771").Value = 86

# Drop the temporary "Text" number format we applied above so the cell keeps
# its plain (unstyled) text value, matching the rest of the sheet.
$ws.Range("A771").ClearFormats()
